$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Receipt ID (A2) and Reference code (O2) values
$ws.Range("A2").Value = "AB522581043"
$ws.Range("O2").Value = "tester16"

# Update the selected cell to B2
$ws.Range("B2").Select()
